$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D1").Value = "Nome"

# Updated data for rows 2-9 (columns C and D)
$data = @(
    @{ Row = 2;  C = "bamboo toilet paper 5 ply 50m bamboo core 100 percent bamboo pulp plastic free FSC Ecolabel OEM"; D = "Carta igienica" },
    @{ Row = 3;  C = "bamboo jumbo tissue roll large and mini jumbo 100 percent bamboo pulp plastic free FSC OEM"; D = "Rotolo jumbo" },
    @{ Row = 4;  C = "bamboo paper hand towels roll or folded 100 percent bamboo pulp plastic free FSC OEM"; D = "Asciugamani carta" },
    @{ Row = 5;  C = "A4 copy paper 80gsm 100 percent bamboo pulp plastic free FSC Ecolabel OEM custom logo"; D = "Carta A4" },
    @{ Row = 6;  C = "notebooks and bloc-notes bamboo paper pages kraft cover plastic free FSC custom logo"; D = "Quaderni blocchi" },
    @{ Row = 7;  C = "facial tissues 100 percent bamboo pulp pocket or box plastic free FSC Ecolabel OEM"; D = "Fazzoletti naso" },
    @{ Row = 8;  C = "kraft paper tape water-activated gummed biodegradable plastic free FSC custom print"; D = "Nastro kraft" },
    @{ Row = 9;  C = "bamboo kraft recycled paper packaging boxes and mailers plastic free FSC custom branding"; D = "Packaging carta" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}

# Remove rows 10-12 (now obsolete entries), deleting entire rows
$ws.Range("A10:D12").EntireRow.Delete()
